$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & 1h volume change) scraped on
# Tue Nov 28 03:53:31 UTC 2023. Column D = Price, Column E = Volume(1h).
# A handful of rows (15-17) also rotated coin/name/link because the
# underlying ranking order changed.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.122.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.024.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.314.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.00%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.742"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.17%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.020.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.069.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0814"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.37"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.12%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.09%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.50"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0609"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.468.16"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.80"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.38"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.29%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.20%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.67"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.209.16"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.02%  "

